$d = $word.ActiveDocument

# 1) Remove "meget sent i projektet, så " between "udført " and "da den blev udført..."
$d.Content.Find.Execute(
    "Forbrugertesten blev udført meget sent i projektet, så da den blev udført, var det allerede besluttet, at udviklingen",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Forbrugertesten blev udført da udviklingen",
    2
)

# 2) Turn "var påbegyndt, og havde højeste prioritet" into "var af højeste prioritet"
$d.Content.Find.Execute(
    "var påbegyndt, og havde højeste prioritet",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "var af højeste prioritet",
    2
)
